$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Swap B32 and C32 values (Non/Oui <-> Oui/Non)
$bVal = $ws.Range("B32").Value2
$cVal = $ws.Range("C32").Value2
$ws.Range("B32").Value = $cVal
$ws.Range("C32").Value = $bVal

# Add new remark text in E32
$ws.Range("E32").Value = "Modifié a 15h (03/07/15)"

# Adjust column E width to (best-fit) width matching the new content
$ws.Columns.Item(5).ColumnWidth = 22

# Update the view: scroll position + selection
$ws.Range("B20").Select()
$ws.Range("D37").Select()
